$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "Tyrael" -> "Tyreal" for every card whose expansion/owner
# column (F) references that name (rows 2:25). Re-assigning every
# occurrence lets the shared-strings table drop the old "Tyrael" entry
# entirely (nothing else references it) and append the corrected
# "Tyreal" string, which is exactly the shared-string churn in the diff.
$ws.Range("F2:F25").Value = "Tyreal"

# "fixed decimals on AC": explicitly (re)apply the General number format
# to the AC value in D4 so it no longer inherits/shows spurious decimals.
$ws.Range("D4").NumberFormat = "General"

# Restore the view: no frozen/scrolled topLeftCell, and the current
# selection sitting on L8.
$ws.Range("L8").Select()
